$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the purpose column (E) and libraryPreparer column (B) for data rows 2-24.
# Order matters for shared-string table layout: "fullRNASEQ" must be interned
# before "S.GISH" so new shared-string indices come out as 30 and 31 respectively.
for ($r = 2; $r -le 24; $r++) {
    $ws.Cells.Item($r, 5).Value = "fullRNASEQ"
    $ws.Cells.Item($r, 2).Value = "S.GISH"
}

$ws.Range("B3:B24").Select()
